$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("weibull")
$ws1.Range("B2").Value = -2.33006830871177
$ws1.Range("C2").Value = 0.10533870155035
$ws1.Range("B3").Value = 0.105924771630349
$ws1.Range("C3").Value = 0.0881247958125287

$ws2 = $wb.Worksheets.Item("lognormal")
$ws2.Range("B2").Value = 1.97970353858253
$ws2.Range("C2").Value = 0.145690465580844
$ws2.Range("B3").Value = -1.12832256940951
$ws2.Range("C3").Value = 0.103194604784936

$ws3 = $wb.Worksheets.Item("llogis")
$ws3.Range("B2").Value = -1.72519678835276
$ws3.Range("C2").Value = 0.0869111420140341
$ws3.Range("B3").Value = 1.99546287416998
$ws3.Range("C3").Value = 0.204557243415721

$ws4 = $wb.Worksheets.Item("gompertz")
$ws4.Range("B2").Value = -2.01765690253255
$ws4.Range("C2").Value = 0.103447449067094
$ws4.Range("B3").Value = -0.0190740131440861
$ws4.Range("C3").Value = 0.0138713780231997

$ws6 = $wb.Worksheets.Item("weibull cov")
$ws6.Range("A2").Value = 0.0110962420443137
$ws6.Range("B2").Value = -0.00350576567788289
$ws6.Range("A3").Value = -0.00350576567788289
$ws6.Range("B3").Value = 0.00776597963699988

$ws7 = $wb.Worksheets.Item("lognormal cov")
$ws7.Range("A2").Value = 0.0212257117611632
$ws7.Range("B2").Value = -0.0124958410297487
$ws7.Range("A3").Value = -0.0124958410297487
$ws7.Range("B3").Value = 0.0106491264567192

$ws8 = $wb.Worksheets.Item("llogis cov")
$ws8.Range("A2").Value = 0.00755354660618361
$ws8.Range("B2").Value = 0.00949575258569563
$ws8.Range("A3").Value = 0.00949575258569563
$ws8.Range("B3").Value = 0.0418436658338387

$ws9 = $wb.Worksheets.Item("gompertz cov")
$ws9.Range("A2").Value = 0.0107013747184891
$ws9.Range("B2").Value = -0.00050452431854454
$ws9.Range("A3").Value = -0.00050452431854454
$ws9.Range("B3").Value = 0.000192415128262507

$wb.Save()
